$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Final data for rows 2..24: item_number, video_file, correct, type, difficulty ---
# (8 new stimuli added: 2 anger, 2 fear, 2 joy, 2 sadness; practice item_numbers renumbered)
$rows = @(
    @{ Row=2;  A=1;   B="adagio04Angry";  C="anger";   D="test";     E=20 }
    @{ Row=3;  A=2;   B="adagio10Angry";  C="anger";   D="test";     E=36.67 }
    @{ Row=4;  A=3;   B="adagio02Angry";  C="anger";   D="test";     E=46.67 }
    @{ Row=5;  A=4;   B="adagio11Angry2"; C="anger";   D="test";     E=30 }
    @{ Row=6;  A=5;   B="adagio13Angry";  C="anger";   D="test";     E=36.67 }
    @{ Row=7;  A=6;   B="adagio20Fear2";  C="fear";    D="test";     E=23.33 }
    @{ Row=8;  A=7;   B="adagio07Fear";   C="fear";    D="test";     E=43.33 }
    @{ Row=9;  A=8;   B="adagio16Fear";   C="fear";    D="test";     E=50 }
    @{ Row=10; A=9;   B="adagio22Fear";   C="fear";    D="test";     E=40 }
    @{ Row=11; A=10;  B="adagio13Fear";   C="fear";    D="test";     E=43.33 }
    @{ Row=12; A=11;  B="adagio24Happy";  C="joy";     D="test";     E=16.670000000000002 }
    @{ Row=13; A=12;  B="adagio13Happy2"; C="joy";     D="test";     E=40 }
    @{ Row=14; A=13;  B="adagio22Happy";  C="joy";     D="test";     E=50 }
    @{ Row=15; A=14;  B="adagio07Happy";  C="joy";     D="test";     E=40 }
    @{ Row=16; A=15;  B="adagio20Happy";  C="joy";     D="test";     E=43.33 }
    @{ Row=17; A=16;  B="adagio18Sad";    C="sadness"; D="test";     E=16.670000000000002 }
    @{ Row=18; A=17;  B="adagio07Sad";    C="sadness"; D="test";     E=40 }
    @{ Row=19; A=18;  B="adagio21Sad";    C="sadness"; D="test";     E=46.67 }
    @{ Row=20; A=19;  B="adagio05Sad2";   C=$null;     D="test";     E=36.67 }
    @{ Row=21; A=20;  B="adagio17Sad";    C=$null;     D="test";     E=40 }
    @{ Row=22; A=100; B="adagio18Angry";  C="anger";   D="practice"; E=$null }
    @{ Row=23; A=101; B="adagio11Happy";  C="joy";     D="practice"; E=$null }
    @{ Row=24; A=102; B="adagio04Sad";    C="sadness"; D="practice"; E=$null }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # A: item_number - literal for row 2 and rows 22-24 (practice); formula elsewhere
    if ($rowNum -eq 2) {
        $ws.Range("A2").Value = $r.A
    } elseif ($rowNum -eq 3) {
        $ws.Range("A3").Formula = "=A2+1"
    } elseif ($rowNum -ge 4 -and $rowNum -le 21) {
        # handled below in one shot (shared formula)
    } else {
        $ws.Cells.Item($rowNum, 1).Value = $r.A
    }

    $ws.Cells.Item($rowNum, 2).Value = $r.B

    if ($null -eq $r.C) {
        $ws.Cells.Item($rowNum, 3).ClearContents()
    } else {
        $ws.Cells.Item($rowNum, 3).Value = $r.C
    }

    $ws.Cells.Item($rowNum, 4).Value = $r.D

    if ($null -eq $r.E) {
        $ws.Cells.Item($rowNum, 5).ClearContents()
    } else {
        $ws.Cells.Item($rowNum, 5).Value = $r.E
    }
}

# A4:A21 share one relative formula "=A3+1" (Excel collapses this into a shared formula group)
$ws.Range("A4:A21").Formula = "=A3+1"

# Column B no longer has a custom width override (back to sheet default)
$ws.Columns("B").UseStandardWidth = $true

# Selection moves to F11 (single cell, was F2:F13 with active cell F13)
$ws.Range("F11").Select()
